$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1025
    $ws.Range("F3").Value = 10
    $ws.Range("F4").Value = 497
    $ws.Range("F5").Value = 439
}
